# "Miglioramenti e pulizia generale"
# - Clear the leftover "${String}" placeholder text out of cell A3
#   (sharedStrings collapses from 2 -> 1 unique string as a result).
# - Move the active selection from A3 to G10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").ClearContents()

$ws.Range("G10").Select()
